$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 0.5604026478897502
$ws.Range("D2").Value = 0.5808675041729621

# Row 3
$ws.Range("C3").Value = -0.4112874021182295
$ws.Range("D3").Value = 0.6848408353172135

# Row 4
$ws.Range("C4").Value = 0.6951750219158087
$ws.Range("D4").Value = 0.4942235038633283

# Row 5
$ws.Range("C5").Value = -1.773796356016656
$ws.Range("D5").Value = 0.08994071619129151

# Row 6
$ws.Range("C6").Value = -0.8062310368019869
$ws.Range("D6").Value = 0.4287391610151932

# Row 7
$ws.Range("C7").Value = 0.1686583835982415
$ws.Range("D7").Value = 0.8676063775620082

# Row 8
$ws.Range("C8").Value = -1.99371521027555
$ws.Range("D8").Value = 0.05872579719085746
$ws.Range("G8").Value = "No"

# Row 9
$ws.Range("C9").Value = 0.7810903848870433
$ws.Range("D9").Value = 0.4430739304220652

# Row 10
$ws.Range("C10").Value = -1.560469513570037
$ws.Range("D10").Value = 0.1329197865475928

# Row 11
$ws.Range("C11").Value = -2.003904972693956
$ws.Range("D11").Value = 0.05754613733241176
$ws.Range("G11").Value = "No"
